$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# "Generate Report for Handoff"
#
# Two en-US source docs were re-handed-off under new GUIDs since the last
# report ran:
#   59e6fecf-c189-4854-80db-5915e8ed3e9f.md -> e6f30ad4-923a-4500-9d07-14030507f370.md
#   89e7c4ae-850b-43df-9d17-84df44657eb4.md -> ffff95e8289a-0bb2-4d81-a22f-748725d7e54d.md
# and their combined xlf package is now:
#   e6f30ad4-923a-4500-9d07-14030507f370.bd19de4e2d5fad4ee8228e061147df4a76f02433.{lang}.xlf
# Status moves from "Handed back: in sync with en-US" to "Ready for handoff",
# and since the new package hasn't been handed back yet, the per-language
# "Latest Target File" / "Latest Handback File" columns are cleared and the
# handback timestamp resets to the zero date.
# -------------------------------------------------------------------------

$mdNameA = "e6f30ad4-923a-4500-9d07-14030507f370.md"
$mdNameB = "ffff95e8289a-0bb2-4d81-a22f-748725d7e54d.md"
$xlfBase = "e6f30ad4-923a-4500-9d07-14030507f370.bd19de4e2d5fad4ee8228e061147df4a76f02433"

$statusText   = "Ready for handoff"
$handoffDate  = "2016-03-24 23:15:16"
$zhHandoffDt  = "2016-03-24 23:15:09"
$zeroDateTime = "0001-01-01 00:00:00"

$mdCommit = "48847690db1acef75c4a86212bb0eebf366561bd"

# =========================================================================
# Sheet "Overview"
# =========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B2").Value = $statusText
$ws1.Range("C2").Value = $statusText
$ws1.Range("D2").Value = $handoffDate
$ws1.Range("B3").Value = $statusText
$ws1.Range("C3").Value = $statusText
$ws1.Range("D3").Value = $handoffDate

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$mdNameA", "", "", $mdNameA)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$mdNameB", "", "", $mdNameB)

# =========================================================================
# Sheet "zh-cn"
# =========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$zhXlfName = "$xlfBase.zh-cn.xlf"
$zhCommit  = "85ab316bc0d3e0c5cc6cb5e16989936acc024a5a"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"

$ws2.Range("C2").Value = $statusText
$ws2.Range("C3").Value = $statusText

$ws2.Range("D2").Value = $zhXlfName
$ws2.Range("D3").Value = $zhXlfName

$ws2.Range("E2").Value = $zhHandoffDt
$ws2.Range("E3").Value = $zhHandoffDt

$ws2.Range("H2").Value = $zeroDateTime
$ws2.Range("H3").Value = $zeroDateTime

# "Latest Target File" (F) / "Latest Handback File" (G) are no longer
# populated -- clear style back to default so the cells drop out entirely
# rather than lingering as empty-but-styled cells.
$ws2.Range("F2:G3").Style = "Normal"
$ws2.Range("F2:G3").ClearContents()

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$mdNameA", "", "", $mdNameA)
$ws2.Hyperlinks.Add($ws2.Range("D2"), $zhXlfUrl, "", "", $zhXlfName)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$mdNameB", "", "", $mdNameB)
$ws2.Hyperlinks.Add($ws2.Range("D3"), $zhXlfUrl, "", "", $zhXlfName)

# =========================================================================
# Sheet "de-de"
# =========================================================================
$ws3 = $wb.Worksheets.Item("de-de")

$deXlfName = "$xlfBase.de-de.xlf"
$deCommit  = "32fc85e0f14443579c6de2a8f3c0c96d990f58d3"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$ws3.Range("C2").Value = $statusText
$ws3.Range("C3").Value = $statusText

$ws3.Range("D2").Value = $deXlfName
$ws3.Range("D3").Value = $deXlfName

# "Latest Handoff Datetime" on this sheet shares the same underlying value
# as the Overview sheet's "Latest Handoff Date" column (both were
# 2016-03-24 23:12:52 before, both become 2016-03-24 23:15:16 now).
$ws3.Range("E2").Value = $handoffDate
$ws3.Range("E3").Value = $handoffDate

$ws3.Range("H2").Value = $zeroDateTime
$ws3.Range("H3").Value = $zeroDateTime

$ws3.Range("F2:G3").Style = "Normal"
$ws3.Range("F2:G3").ClearContents()

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$mdNameA", "", "", $mdNameA)
$ws3.Hyperlinks.Add($ws3.Range("D2"), $deXlfUrl, "", "", $deXlfName)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$mdNameB", "", "", $mdNameB)
$ws3.Hyperlinks.Add($ws3.Range("D3"), $deXlfUrl, "", "", $deXlfName)
